$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C90").Value = 36956

$ws.Range("B91").Value = 41600

$ws.Range("A92").Value = 43897
$ws.Range("B92").Value = 18720
$ws.Range("D92").Value = $ws.Range("D88").Value
$ws.Range("E92").Formula = $ws.Range("E91").Formula

$ws.Range("A93").Value = 43897
$ws.Range("C93").Value = 36956
$ws.Range("D93").Value = $ws.Range("D90").Value
$ws.Range("E93").Formula = $ws.Range("E91").Formula
